$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "保險" (insurance) -- add company/name/owner/property_category/
# category/date/legislator_name/legislator_id/source_file/index columns
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("保險")

$headers5 = @("company","name","owner","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headers5.Length; $i++) {
    $ws5.Cells.Item(1, 2 + $i).Value = $headers5[$i]
}

$rows5 = @(
    @{ company = "三商美邦人壽"; name = "世紀理財變額萬能終身壽險"; owner = "蔡煌瑯"; property_category = "insurance"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 80 },
    @{ company = "三商美邦人壽"; name = "世紀理財變額萬能終身壽險"; owner = "王琴賀"; property_category = "insurance"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 81 },
    @{ company = "新光人壽";     name = "美利外幣終生還本型保險"; owner = "王琴賀"; property_category = "insurance"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 82 }
)

for ($r = 0; $r -lt $rows5.Length; $r++) {
    $row = $r + 2
    $data = $rows5[$r]
    for ($i = 0; $i -lt $headers5.Length; $i++) {
        $ws5.Cells.Item($row, 2 + $i).Value = $data[$headers5[$i]]
    }
}

# ---------------------------------------------------------------------------
# Sheet "債務" (debt) -- add species/debtor/owner/total/register_date/
# register_reason/property_category/category/date/legislator_name/
# legislator_id/source_file/index columns
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("債務")

$headers6 = @("species","debtor","owner","total","register_date","register_reason","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $headers6.Length; $i++) {
    $ws6.Cells.Item(1, 2 + $i).Value = $headers6[$i]
}

$rows6 = @(
    @{ species = "房屋貸款"; debtor = "王琴賀"; owner = "台中商業銀行埔里分行南投縣埔里鎮西康路62號"; total = 3644844;  register_date = "95年01月23日";  register_reason = "設定"; property_category = "debt"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 92 },
    @{ species = "房屋貸款"; debtor = "王琴賀"; owner = "台中商業銀行埔里分行南投縣埔里鎮西康路62號"; total = 307104;   register_date = "98年11月16日";  register_reason = "設定"; property_category = "debt"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 93 },
    @{ species = "房屋貸款"; debtor = "蔡煌瑯"; owner = "臺灣銀行臺北市中正區重慶南路一段120號";         total = 13056565; register_date = "98年01月16日";  register_reason = "設定"; property_category = "debt"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 94 },
    @{ species = "房屋貸款"; debtor = "王琴賀"; owner = "台中商業銀行埔里分行南投縣埔里鎮西康路62號"; total = 850002;   register_date = "99年01月15日";  register_reason = "設定"; property_category = "debt"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 95 },
    @{ species = "房屋貸款"; debtor = "蔡煌瑯"; owner = "臺灣銀行臺北市中正區重慶南路一段120號";         total = 2918530;  register_date = "102年03月28日"; register_reason = "設定"; property_category = "debt"; category = "normal"; date = "2013-12-17"; legislator_name = "蔡煌瑯"; legislator_id = 752; source_file = "tmpc9fc1"; index = 97 }
)

for ($r = 0; $r -lt $rows6.Length; $r++) {
    $row = $r + 2
    $data = $rows6[$r]
    for ($i = 0; $i -lt $headers6.Length; $i++) {
        $ws6.Cells.Item($row, 2 + $i).Value = $data[$headers6[$i]]
    }
}

Write-Host "Edit complete"
